$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '34.164.06'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.792.04'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.72%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '224.73'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.53%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.38%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '32.71'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.29%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.286'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0708'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.46%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0930'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.07%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.048.58'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.57%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.799.17'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.96%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.87'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -2.04%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.625'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -2.99%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '34.134.99'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.25%  '
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.81%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '68.06'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.66%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '243.23'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -3.39%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0785'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.77%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.05%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.73'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -4.09%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.10'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -3.44%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -2.84%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '159.59'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '16.28'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.79%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.06'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.35%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.113'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -1.44%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.14%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0517'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -1.74%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +1.52%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.68'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -3.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.51'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -2.77%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.82'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -4.01%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -2.19%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.648'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.66%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -1.43%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0187'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.30%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.21'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +3.09%  '
$ws.Range('B40').Value = 'HuobiToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.34'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.45%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '79.04'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -3.70%  '
$ws.Range('B42').Value = 'ARBITRUM'
$ws.Range('C42').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.917'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -5.04%  '
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.70'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -3.01%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0₆0146'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +17.90%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.86%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '108.20'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +2.71%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0496'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.39%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.88'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -3.04%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.948.17'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.05%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '12.17'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.20%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.05%  '
